$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new value.
$updates = @(
    @("D2", "42.976.04"),
    @("E2", "  +2.00%  "),
    @("D3", "2.299.19"),
    @("E3", "  +1.41%  "),
    @("E4", "  +0.10%  "),
    @("D5", "309.04"),
    @("E5", "  +1.21%  "),
    @("D6", "100.54"),
    @("E6", "  +4.54%  "),
    @("D7", "0.533"),
    @("E7", "  +0.61%  "),
    @("E8", "  +0.00%  "),
    @("D9", "0.506"),
    @("E9", "  +3.28%  "),
    @("D10", "35.83"),
    @("E10", "  +1.67%  "),
    @("D11", "0.0818"),
    @("E11", "  +2.97%  "),
    @("E12", "  +0.52%  "),
    @("D13", "6.95"),
    @("E13", "  +4.66%  "),
    @("D14", "2.659.45"),
    @("E14", "  +1.58%  "),
    @("D15", "14.83"),
    @("E15", "  +3.01%  "),
    @("D16", "2.299.07"),
    @("E16", "  +0.57%  "),
    @("D17", "0.801"),
    @("E17", "  +0.98%  "),
    @("D18", "42.973.06"),
    @("E18", "  +2.16%  "),
    @("D19", "12.51"),
    @("E19", "  +0.47%  "),
    @("E21", "  +1.51%  "),
    @("D22", "68.04"),
    @("E22", "  +0.34%  "),
    @("D23", "239.45"),
    @("E23", "  +0.55%  "),
    @("D24", "2.01"),
    @("E24", "  +3.86%  "),
    @("D25", "2.60"),
    @("E25", "  +0.98%  "),
    @("D26", "0.992"),
    @("E26", "  -0.66%  "),
    @("D27", "24.26"),
    @("E27", "  +2.35%  "),
    @("D28", "38.81"),
    @("E28", "  +4.80%  "),
    @("D29", "9.62"),
    @("E29", "  +1.14%  "),
    @("D30", "2.11"),
    @("E30", "  +0.01%  "),
    @("D31", "165.22"),
    @("E31", "  +3.70%  "),
    @("D32", "5.31"),
    @("E32", "  +1.12%  "),
    @("E33", "  +0.23%  "),
    @("D34", "3.15"),
    @("E34", "  -1.39%  "),
    @("D35", "17.69"),
    @("E35", "  +3.45%  "),
    @("D36", "0.0738"),
    @("E36", "  -0.27%  "),
    @("D37", "2.39"),
    @("E38", "  -0.04%  "),
    @("D39", "1.84"),
    @("E39", "  +0.76%  "),
    @("E40", "  +1.05%  "),
    @("E41", "  +3.94%  "),
    @("E42", "  -5.71%  "),
    @("D43", "0.0289"),
    @("E43", "  +2.06%  "),
    @("D44", "1.966.57"),
    @("E44", "  -1.12%  "),
    @("E45", "  +1.42%  "),
    @("D46", "3.00"),
    @("E46", "  +2.90%  "),
    @("D47", "9.79"),
    @("E47", "  -1.45%  "),
    @("D48", "2.99"),
    @("E48", "  +18.93%  "),
    @("D49", "54.88"),
    @("E49", "  +3.23%  "),
    @("B50", "RocketPoolETH"),
    @("C50", "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"),
    @("D50", "2.526.28"),
    @("E50", "  +1.56%  "),
    @("B51", "Stacks"),
    @("C51", "https://coinranking.com/coin/mMPrMcB7+stacks-stx"),
    @("D51", "1.54"),
    @("E51", "  +1.80%  ")
)

foreach ($u in $updates) {
    $cell = $u[0]
    $value = $u[1]
    $range = $ws.Range($cell)
    if ($cell.Substring(0,1) -eq "D") {
        # Column D holds price text that can look numeric (e.g. "309.04" or "3.00").
        # Force text storage so precision/trailing zeros survive, same as the source inline string,
        # then restore the default "Normal" style so no stray number format sticks to the cell.
        $range.NumberFormat = "@"
        $range.Value = $value
        $range.Style = "Normal"
    } else {
        $range.Value = $value
    }
}
